# Apply the "hourly rate conditional template" edit:
# Replace the literal field {{ person.hourly_rate }} in the evaluator table
# with a conditional expression that only prints a formatted currency value
# when person.hourly_rate is truthy:
#   {% if person.hourly_rate%}{{ currency (person.hourly_rate) }}{% endif %}

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "{{ person.hourly_rate }}",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "{% if person.hourly_rate%}{{ currency (person.hourly_rate) }}{% endif %}",
    2
)
